$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated TPM-derived values for all Sending-cluster x Target-cluster combinations,
# plus a new "Resolving-Mac" target/sending cluster (adds rows 11-13, shared string 25).
$rows = @(
  @("ECs", "Efna5", "Epha3", "ECs", 1, 0.3333333333333333, 0.1728506666666667, 0.518552, 0.0840503369699626, 0.0840503369699626, 2, 0.6666666666666666, 0.07579599999999999, 0.227388, 0.001780200955210419, 0.001780200955210419, 0.01310138913066667, 0.117912502176, 0.000149626490159685, 0.000149626490159685),
  @("ECs", "Efna5", "Epha3", "FAPs", 1, 0.3333333333333333, 0.1728506666666667, 0.518552, 0.0840503369699626, 0.0840503369699626, 3, 1, 40.12734033333333, 120.382021, 0.9424604146848589, 0.9424604146848587, 6.936037528176889, 62.424337753592, 0.07921411543511307, 0.07921411543511307),
  @("ECs", "Efna5", "Epha3", "MuSCs", 1, 0.3333333333333333, 0.1728506666666667, 0.518552, 0.0840503369699626, 0.0840503369699626, 3, 1, 2.344072666666666, 7.032217999999999, 0.05505462557763778, 0.05505462557763778, 0.4051745231484444, 3.646570708336, 0.004627359831555577, 0.004627359831555577),
  @("ECs", "Efna5", "Epha3", "Resolving-Mac", 1, 0.3333333333333333, 0.1728506666666667, 0.518552, 0.0840503369699626, 0.0840503369699626, 2, 0.6666666666666666, 0.03000666666666667, 0.09002, 0.0007047587822930054, 0.0007047587822930053, 0.005186672337777779, 0.04668005104, 0.00005923521313426762, 0.00005923521313426761),
  @("FAPs", "Efna5", "Epha3", "ECs", 3, 1, 1.367717666666667, 4.103153, 0.6650661694281633, 0.6650661694281633, 2, 0.6666666666666666, 0.07579599999999999, 0.227388, 0.001780200955210419, 0.001780200955210419, 0.1036675282626666, 0.9330077543639999, 0.001183951430094151, 0.001183951430094151),
  @("FAPs", "Efna5", "Epha3", "FAPs", 3, 1, 1.367717666666667, 4.103153, 0.6650661694281633, 0.6650661694281633, 3, 1, 40.12734033333333, 120.382021, 0.9424604146848589, 0.9424604146848587, 54.88287229024588, 493.945850612213, 0.6267985378321375, 0.6267985378321373),
  @("FAPs", "Efna5", "Epha3", "MuSCs", 3, 1, 1.367717666666667, 4.103153, 0.6650661694281633, 0.6650661694281633, 3, 1, 2.344072666666666, 7.032217999999999, 0.05505462557763778, 0.05505462557763778, 3.206029598150444, 28.85426638335399, 0.03661496894222134, 0.03661496894222134),
  @("FAPs", "Efna5", "Epha3", "Resolving-Mac", 3, 1, 1.367717666666667, 4.103153, 0.6650661694281633, 0.6650661694281633, 2, 0.6666666666666666, 0.03000666666666667, 0.09002, 0.0007047587822930054, 0.0007047587822930053, 0.04104064811777778, 0.36936583306, 0.000468711223710466, 0.000468711223710466),
  @("MuSCs", "Efna5", "Epha3", "ECs", 3, 1, 0.5159453333333334, 1.547836, 0.2508834936018741, 0.2508834936018741, 2, 0.6666666666666666, 0.07579599999999999, 0.227388, 0.001780200955210419, 0.001780200955210419, 0.03910659248533333, 0.351959332368, 0.0004466230349565833, 0.0004466230349565833),
  @("MuSCs", "Efna5", "Epha3", "FAPs", 3, 1, 0.5159453333333334, 1.547836, 0.2508834936018741, 0.2508834936018741, 3, 1, 40.12734033333333, 120.382021, 0.9424604146848589, 0.9424604146848587, 20.70351398406178, 186.331625856556, 0.2364477614176084, 0.2364477614176084),
  @("MuSCs", "Efna5", "Epha3", "MuSCs", 3, 1, 0.5159453333333334, 1.547836, 0.2508834936018741, 0.2508834936018741, 3, 1, 2.344072666666666, 7.032217999999999, 0.05505462557763778, 0.05505462557763778, 1.209413353360889, 10.884720180248, 0.01381229680386086, 0.01381229680386086),
  @("MuSCs", "Efna5", "Epha3", "Resolving-Mac", 3, 1, 0.5159453333333334, 1.547836, 0.2508834936018741, 0.2508834936018741, 2, 0.6666666666666666, 0.03000666666666667, 0.09002, 0.0007047587822930054, 0.0007047587822930053, 0.01548179963555556, 0.13933619672, 0.0001768123454482718, 0.0001768123454482718)
)

for ($i = 0; $i -lt $rows.Count; $i++) {
    $row = $rows[$i]
    $r = $i + 2
    for ($j = 0; $j -lt $row.Count; $j++) {
        $ws.Cells.Item($r, $j + 1).Value = $row[$j]
    }
}
